# Commit: "ajuste pesp y firma documentos codigo ado digicredito"
# Updates the single data row (row 2) in sheet "dataDigicredito" with a new
# set of applicant/loan test-data values, and adjusts the active window
# selection/column widths used while reviewing that data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataDigicredito")

# --- Row 2 data updates -----------------------------------------------
$ws.Range("B2").Value  = '"3994518"'              # Cedula
$ws.Range("D2").Value  = '"16"'                   # Plazo
$ws.Range("F2").Value  = '"90"'                   # DiasHabilesIntereses
$ws.Range("L2").Value  = '"SERGIO"'                # primerNombre
$ws.Range("M2").Value  = '"SEGUNDO"'               # segundoNombre
$ws.Range("N2").Value  = '"URZOLA"'                # primerApellido
$ws.Range("O2").Value  = '"BERTEL"'                # segundoApellido
$ws.Range("P2").Value  = '"18/Mar/2022"'           # fechaActual
$ws.Range("T2").Value  = '"500000"'                # vlrCompasSaneamientos
$ws.Range("W2").Value  = '"lespitiameza@gmail.com" ' # Correo
$ws.Range("X2").Value  = '"3219176522"'            # Celular
$ws.Range("AX2").Value = '"Compra de cartera"'     # lineaCredito
$ws.Range("CW2").Value = '"Si"'                    # tomarSeguroAP

# --- Column width tweaks (best fit on Cedula / Credito / fechaActual) --
$ws.Columns.Item(1).ColumnWidth = 19.7109375
$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(16).ColumnWidth = 13.85546875

# --- Active window selection -------------------------------------------
$ws.Activate()
$ws.Range("CW9").Select()
